# Add copies option, ref other options
# Append new usage-log rows (52-67) to the "Наличные" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Наличные")

$data = @(
    @(52, 7676096317, "6, 7.docx", 4, 0.8, "2025-07-07 01:17:53"),
    @(53, 7676096317, "6, 7.docx", 4, 0.8, "2025-07-07 01:18:42"),
    @(54, 7676096317, "6, 7.docx", 4, 0.8, "2025-07-07 01:19:40"),
    @(55, 7676096317, "6, 7.docx", 4, 0.8, "2025-07-07 01:25:15"),
    @(56, 7676096317, "3.docx", 2, 0.4, "2025-07-07 01:46:42"),
    @(57, 7676096317, "6, 7.docx", 4, 0.8, "2025-07-07 01:47:25"),
    @(58, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:07:11"),
    @(59, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:09:26"),
    @(60, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:11:25"),
    @(61, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:34:15"),
    @(62, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:47:47"),
    @(63, 7676096317, "debug.pdf", 4, 0.8, "2025-07-09 15:49:00"),
    @(64, 7676096317, "Radkovich_Otchetik.docx", 33, 6.6, "2025-07-09 15:51:28"),
    @(65, 7676096317, "Radkovich_Otchetik.docx", 33, 6.6, "2025-07-09 15:54:07"),
    @(66, 7676096317, "Radkovich_Otchetik.docx", 33, 6.6, "2025-07-09 15:54:47"),
    @(67, 7676096317, "Radkovich_Otchetik.docx", 33, 6.6, "2025-07-09 15:57:15")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
